$d = $word.ActiveDocument

$d.Content.Find.Execute("380÷7=54, 2", $true, $false, $false, $false, $false, $true, 1, $false, "561÷4=140, 1", 2) | Out-Null
$d.Content.Find.Execute("538÷9=59, 7", $true, $false, $false, $false, $false, $true, 1, $false, "900÷8=112, 4", 2) | Out-Null
$d.Content.Find.Execute("179÷6=29, 5", $true, $false, $false, $false, $false, $true, 1, $false, "798÷8=99, 6", 2) | Out-Null
$d.Content.Find.Execute("866÷8=108, 2", $true, $false, $false, $false, $false, $true, 1, $false, "852÷6=142, 0", 2) | Out-Null
$d.Content.Find.Execute("969÷7=138, 3", $true, $false, $false, $false, $false, $true, 1, $false, "543÷2=271, 1", 2) | Out-Null
$d.Content.Find.Execute("840÷6=140, 0", $true, $false, $false, $false, $false, $true, 1, $false, "213÷5=42, 3", 2) | Out-Null
$d.Content.Find.Execute("315÷8=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "277÷7=39, 4", 2) | Out-Null
$d.Content.Find.Execute("160÷7=22, 6", $true, $false, $false, $false, $false, $true, 1, $false, "782÷4=195, 2", 2) | Out-Null
$d.Content.Find.Execute("221÷8=27, 5", $true, $false, $false, $false, $false, $true, 1, $false, "355÷9=39, 4", 2) | Out-Null
$d.Content.Find.Execute("953÷9=105, 8", $true, $false, $false, $false, $false, $true, 1, $false, "769÷4=192, 1", 2) | Out-Null
$d.Content.Find.Execute("180÷5=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "945÷7=135, 0", 2) | Out-Null
$d.Content.Find.Execute("179÷9=19, 8", $true, $false, $false, $false, $false, $true, 1, $false, "377÷4=94, 1", 2) | Out-Null
$d.Content.Find.Execute("306÷7=43, 5", $true, $false, $false, $false, $false, $true, 1, $false, "146÷6=24, 2", 2) | Out-Null
$d.Content.Find.Execute("597÷5=119, 2", $true, $false, $false, $false, $false, $true, 1, $false, "424÷6=70, 4", 2) | Out-Null
$d.Content.Find.Execute("846÷4=211, 2", $true, $false, $false, $false, $false, $true, 1, $false, "993÷6=165, 3", 2) | Out-Null
$d.Content.Find.Execute("359÷5=71, 4", $true, $false, $false, $false, $false, $true, 1, $false, "765÷6=127, 3", 2) | Out-Null
$d.Content.Find.Execute("581÷4=145, 1", $true, $false, $false, $false, $false, $true, 1, $false, "312÷5=62, 2", 2) | Out-Null
$d.Content.Find.Execute("840÷8=105, 0", $true, $false, $false, $false, $false, $true, 1, $false, "849÷4=212, 1", 2) | Out-Null
$d.Content.Find.Execute("528÷4=132, 0", $true, $false, $false, $false, $false, $true, 1, $false, "140÷3=46, 2", 2) | Out-Null
$d.Content.Find.Execute("272÷7=38, 6", $true, $false, $false, $false, $false, $true, 1, $false, "588÷2=294, 0", 2) | Out-Null
$d.Content.Find.Execute("141÷6=23, 3", $true, $false, $false, $false, $false, $true, 1, $false, "986÷5=197, 1", 2) | Out-Null
$d.Content.Find.Execute("971÷6=161, 5", $true, $false, $false, $false, $false, $true, 1, $false, "129÷8=16, 1", 2) | Out-Null
$d.Content.Find.Execute("582÷4=145, 2", $true, $false, $false, $false, $false, $true, 1, $false, "553÷7=79, 0", 2) | Out-Null
$d.Content.Find.Execute("666÷4=166, 2", $true, $false, $false, $false, $false, $true, 1, $false, "277÷3=92, 1", 2) | Out-Null
$d.Content.Find.Execute("621÷3=207, 0", $true, $false, $false, $false, $false, $true, 1, $false, "434÷2=217, 0", 2) | Out-Null
